# This script updates the "想去人数" (F column) values on the
# "展览" and "全部类型" worksheets to match a fresh data pull.
$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F2").Value = 30
$wsExhibition.Range("F4").Value = 0
$wsExhibition.Range("F5").Value = 209
$wsExhibition.Range("F6").Value = 0
$wsExhibition.Range("F7").Value = 0
$wsExhibition.Range("F8").Value = 107
$wsExhibition.Range("F9").Value = 92
$wsExhibition.Range("F10").Value = 743
$wsExhibition.Range("F11").Value = 0
$wsExhibition.Range("F12").Value = 1157
$wsExhibition.Range("F14").Value = 0
$wsExhibition.Range("F15").Value = 179
$wsExhibition.Range("F16").Value = 81
$wsExhibition.Range("F17").Value = 146
$wsExhibition.Range("F20").Value = 6263
$wsExhibition.Range("F22").Value = 38
$wsExhibition.Range("F23").Value = 84
$wsExhibition.Range("F24").Value = 0
$wsExhibition.Range("F26").Value = 0
$wsExhibition.Range("F27").Value = 400
$wsExhibition.Range("F29").Value = 9
$wsExhibition.Range("F30").Value = 0
$wsExhibition.Range("F31").Value = 0
$wsExhibition.Range("F32").Value = 528
$wsExhibition.Range("F33").Value = 139
$wsExhibition.Range("F34").Value = 0
$wsExhibition.Range("F36").Value = 369
$wsExhibition.Range("F38").Value = 7
$wsExhibition.Range("F39").Value = 1558
$wsExhibition.Range("F40").Value = 956
$wsExhibition.Range("F44").Value = 0
$wsExhibition.Range("F45").Value = 479
$wsExhibition.Range("F46").Value = 75
$wsExhibition.Range("F47").Value = 581

$wsAllTypes = $wb.Worksheets.Item("全部类型")
$wsAllTypes.Range("F3").Value = 226
$wsAllTypes.Range("F5").Value = 209
$wsAllTypes.Range("F6").Value = 154
$wsAllTypes.Range("F7").Value = 0
$wsAllTypes.Range("F8").Value = 109
$wsAllTypes.Range("F9").Value = 107
$wsAllTypes.Range("F13").Value = 0
$wsAllTypes.Range("F15").Value = 0
$wsAllTypes.Range("F16").Value = 0
$wsAllTypes.Range("F17").Value = 81
$wsAllTypes.Range("F21").Value = 6263
$wsAllTypes.Range("F24").Value = 0
$wsAllTypes.Range("F25").Value = 537
$wsAllTypes.Range("F27").Value = 3961
$wsAllTypes.Range("F28").Value = 0
$wsAllTypes.Range("F34").Value = 0
$wsAllTypes.Range("F35").Value = 280
$wsAllTypes.Range("F37").Value = 0
$wsAllTypes.Range("F38").Value = 169
$wsAllTypes.Range("F39").Value = 0
$wsAllTypes.Range("F40").Value = 1558
$wsAllTypes.Range("F41").Value = 0
$wsAllTypes.Range("F42").Value = 0
$wsAllTypes.Range("F43").Value = 61
$wsAllTypes.Range("F45").Value = 489
$wsAllTypes.Range("F47").Value = 75
$wsAllTypes.Range("F48").Value = 581

